$d = $word.ActiveDocument

# 1. "Above & Beyond:" becomes "Borrowed Code:"
$d.Content.Find.Execute("Above & Beyond:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Borrowed Code:", 2) | Out-Null

# 2. The paragraph that used to read "Borrowed Code:" (the one right after the text we
#    just changed) is removed together with its paragraph mark, so that the following
#    paragraph (the "The parseBody ... " text, with its spell-check run split intact)
#    shifts up to take its place.
$findDup = $d.Content
$findDup.Find.Execute("Borrowed Code:") | Out-Null
$findDup.Collapse(0)
$findDup.Find.Execute("Borrowed Code:") | Out-Null
$oldBorrowed = $findDup.Paragraphs(1)
$d.Range($oldBorrowed.Range.Start, $oldBorrowed.Range.End).Delete() | Out-Null

# 3. Append the new sentence right before the paragraph mark (and the _GoBack bookmark
#    that sits there) of what is now the last paragraph, so it lands in the same run,
#    after the bookmark stays attached to the very end of the paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endOfLastPara = $lastPara.Range.End
$d.Range($endOfLastPara, $endOfLastPara).InsertBefore( `
    "There are other snippets borrowed from Stack Overflow, which are noted in the code comments.")

# 4. Split that paragraph in two right after "...provided in class. " so the new
#    sentence becomes its own, final paragraph (keeping the bookmark at its end).
$splitRange = $d.Content
$splitRange.Find.Execute("class. ") | Out-Null
$splitPos = $splitRange.End
$d.Range($splitPos, $splitPos).InsertParagraphAfter() | Out-Null
